# Updated symbol list (coin rankings/prices/volumes) per the scraper refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text is a plain (non-numeric) string: Coin name / Link columns.
$textUpdates = @(
    @{ Cell = 'B8'; Value = 'MXToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'B10'; Value = 'WazirX' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'B12'; Value = 'BitrueCoin' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'B13'; Value = 'BitMartToken' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'B14'; Value = 'BitForexToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'B15'; Value = 'TigerCash' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'B16'; Value = 'LEO' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'B17'; Value = 'GateToken' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
)

# Cells whose new text looks like a number/percentage (e.g. "291.33", "0.03%").
# Force the cell format to Text ("@") BEFORE assigning so Excel stores the
# exact original string (incl. trailing zeros / % sign) instead of silently
# re-typing the cell as a Number/Percentage and rounding it.
$numericTextUpdates = @(
    @{ Cell = 'D2'; Value = '291.33' }
    @{ Cell = 'E2'; Value = '0.03%' }
    @{ Cell = 'D3'; Value = '30.99' }
    @{ Cell = 'E3'; Value = '0.67%' }
    @{ Cell = 'D4'; Value = '4.956' }
    @{ Cell = 'E4'; Value = '1.18%' }
    @{ Cell = 'D5'; Value = '0.07444' }
    @{ Cell = 'E5'; Value = '2.56%' }
    @{ Cell = 'D6'; Value = '2.238' }
    @{ Cell = 'E6'; Value = '-5.26%' }
    @{ Cell = 'D7'; Value = '7.743' }
    @{ Cell = 'E7'; Value = '0.94%' }
    @{ Cell = 'D8'; Value = '0.9185' }
    @{ Cell = 'E8'; Value = '2.40%' }
    @{ Cell = 'D9'; Value = '0.09367' }
    @{ Cell = 'E9'; Value = '17.43%' }
    @{ Cell = 'D10'; Value = '0.1727' }
    @{ Cell = 'E10'; Value = '3.59%' }
    @{ Cell = 'D11'; Value = '0.08328' }
    @{ Cell = 'E11'; Value = '2.19%' }
    @{ Cell = 'D12'; Value = '0.03209' }
    @{ Cell = 'E12'; Value = '3.74%' }
    @{ Cell = 'D13'; Value = '0.09928' }
    @{ Cell = 'E13'; Value = '-0.90%' }
    @{ Cell = 'D14'; Value = '0.001494' }
    @{ Cell = 'E14'; Value = '-0.55%' }
    @{ Cell = 'D15'; Value = '0.005724' }
    @{ Cell = 'E15'; Value = '-2.03%' }
    @{ Cell = 'D16'; Value = '3.468' }
    @{ Cell = 'E16'; Value = '-0.05%' }
    @{ Cell = 'D17'; Value = '3.762' }
    @{ Cell = 'E17'; Value = '1.59%' }
    @{ Cell = 'D18'; Value = '2.147' }
    @{ Cell = 'E18'; Value = '3.34%' }
    @{ Cell = 'D19'; Value = '0.3331' }
    @{ Cell = 'E19'; Value = '0.38%' }
    @{ Cell = 'D20'; Value = '0.1297' }
    @{ Cell = 'E20'; Value = '-0.08%' }
    @{ Cell = 'D21'; Value = '4.151' }
    @{ Cell = 'E21'; Value = '4.66%' }
    @{ Cell = 'E22'; Value = '-7.98%' }
    @{ Cell = 'D23'; Value = '0.04511' }
    @{ Cell = 'E23'; Value = '0.04%' }
    @{ Cell = 'D24'; Value = '0.001217' }
    @{ Cell = 'E24'; Value = '0.51%' }
    @{ Cell = 'D25'; Value = '0.004256' }
    @{ Cell = 'E25'; Value = '-3.37%' }
    @{ Cell = 'D26'; Value = '0.0001300' }
    @{ Cell = 'D27'; Value = '0.0003388' }
    @{ Cell = 'E27'; Value = '-0.27%' }
    @{ Cell = 'D39'; Value = '0.01620' }
    @{ Cell = 'E39'; Value = '2.49%' }
    @{ Cell = 'D40'; Value = '0.04567' }
    @{ Cell = 'E40'; Value = '4.14%' }
    @{ Cell = 'D41'; Value = '0.007422' }
    @{ Cell = 'E41'; Value = '1.27%' }
    @{ Cell = 'D42'; Value = '0.009821' }
    @{ Cell = 'E42'; Value = '-2.11%' }
    @{ Cell = 'D43'; Value = '0.1356' }
    @{ Cell = 'E43'; Value = '3.25%' }
    @{ Cell = 'D44'; Value = '0.002158' }
    @{ Cell = 'E44'; Value = '6.64%' }
    @{ Cell = 'D45'; Value = '0.008717' }
    @{ Cell = 'E45'; Value = '-8.37%' }
    @{ Cell = 'D46'; Value = '0.00006088' }
    @{ Cell = 'E46'; Value = '6.26%' }
    @{ Cell = 'E47'; Value = '-0.27%' }
    @{ Cell = 'D48'; Value = '2.615' }
    @{ Cell = 'E48'; Value = '16.61%' }
    @{ Cell = 'D49'; Value = '0.001996' }
    @{ Cell = 'E49'; Value = '-31.19%' }
    @{ Cell = 'D50'; Value = '0.00002096' }
    @{ Cell = 'E50'; Value = '-0.27%' }
    @{ Cell = 'D51'; Value = '0.0001996' }
    @{ Cell = 'E51'; Value = '-0.27%' }
)

foreach ($u in $textUpdates) {
    $ws.Range($u.Cell).Value = $u.Value
}

foreach ($u in $numericTextUpdates) {
    $cell = $ws.Range($u.Cell)
    $cell.NumberFormat = "@"
    $cell.Value = $u.Value
}
